# Generate Report for Handoff
# Updates the status/report for the 7a82f988-c0d5-4b60-b7af-c5251f112ba7.md file:
#  - moves it from "Handed back: in sync with en-US" to "Ready for handoff"
#  - refreshes the handoff timestamps
#  - records the stale-handback error detail message
#  - widens the Error Detail column so the new message is readable

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4579062742c0fc5f45923080c2b41405436bb0f6/e2e/7a82f988-c0d5-4b60-b7af-c5251f112ba7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c7aaeb7f4569af33536c8d3d6cc327f4c85e65b/e2e/7a82f988-c0d5-4b60-b7af-c5251f112ba7.md."

# --- Overview sheet: row 3 is the 7a82f988 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 02:58:13"

# --- zh-cn sheet: row 3 is the 7a82f988 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-22 02:58:08"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the 7a82f988 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-22 02:58:13"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
